$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old A1 formula
$ws.Range("A1").ClearContents()

# New data: A2=2, B2=2, C2=A2+B2
$ws.Range("A2").Value = 2
$ws.Range("B2").Value = 2
$ws.Range("C2").Formula = "=A2+B2"

# Update selection to C3 (matches diff's selection activeCell="C3" sqref="C3")
$ws.Range("C3").Select()
